$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 87, pushing the existing rows 87:234 down to 88:235
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with a new data record
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 44665
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112017
$ws.Range("G87").Value = "Apio"
$ws.Range("H87").Value = "Americana (o)"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 25
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 12000
$ws.Range("M87").Value = 12000
$ws.Range("N87").Value = "`$/docena de matas"
$ws.Range("O87").Value = "Región de Coquimbo"
$ws.Range("P87").Value = 2000
$ws.Range("Q87").Value = 6
$ws.Range("R87").Value = "Hortaliza"
